$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the "status" (current) flag from the three EIR Alt1 runs that are
#     no longer the "current" run now that Alt1_03 exists ---
$ws.Range("H135").ClearContents()
$ws.Range("H136").ClearContents()
$ws.Range("H137").ClearContents()

# --- Re-point the previous "current" Alt1 row (previously Alt2_01 data) to
#     the new Alt1_03 run (far-tiers final EIR path) ---
$ws.Range("C138").Value = "2035_TM152_EIR_Alt1_03"
$ws.Range("E138").Value = "Alt1"
$ws.Range("F138").Value = "`"EIR runs\Alt1 (s26) runs\Alt1_v3_test_far_tiers_FINAL_EIR_ALT`""
$ws.Range("G138").Value = "run375"

$ws.Range("C139").Value = "2050_TM152_EIR_Alt1_03"
$ws.Range("E139").Value = "Alt1"
$ws.Range("F139").Value = "`"EIR runs\Alt1 (s26) runs\Alt1_v3_test_far_tiers_FINAL_EIR_ALT`""
$ws.Range("G139").Value = "run375"

# --- Insert four new rows (140-143) below, copying formatting from row 139
#     so the new rows keep the same fill/alignment styling used throughout
#     this block ---
$ws.Rows("139").Copy()
$ws.Rows("140").Insert(-4121, 0)
$ws.Rows("139").Copy()
$ws.Rows("141").Insert(-4121, 0)
$ws.Rows("139").Copy()
$ws.Rows("142").Insert(-4121, 0)
$ws.Rows("139").Copy()
$ws.Rows("143").Insert(-4121, 0)

# Row 140: 2035 EIR Alt2_01 (anticipated run, not yet current)
$ws.Range("A140").Value = "RTP2021"
$ws.Range("B140").Value = 2035
$ws.Range("C140").Value = "2035_TM152_EIR_Alt2_01"
$ws.Range("D140").Value = "EIR"
$ws.Range("E140").Value = "Alt2"
$ws.Range("F140").Value = "`"EIR runs\Alt2 (s28) runs\Alt2_v1`""
$ws.Range("G140").Value = "run374"
$ws.Range("H140").ClearContents()

# Row 141: 2050 EIR Alt2_01 (anticipated run, not yet current)
$ws.Range("A141").Value = "RTP2021"
$ws.Range("B141").Value = 2050
$ws.Range("C141").Value = "2050_TM152_EIR_Alt2_01"
$ws.Range("D141").Value = "EIR"
$ws.Range("E141").Value = "Alt2"
$ws.Range("F141").Value = "`"EIR runs\Alt2 (s28) runs\Alt2_v1`""
$ws.Range("G141").Value = "run374"
$ws.Range("H141").ClearContents()

# Row 142: 2035 EIR Alt2_02 (new anticipated current run)
$ws.Range("A142").Value = "RTP2021"
$ws.Range("B142").Value = 2035
$ws.Range("C142").Value = "2035_TM152_EIR_Alt2_02"
$ws.Range("D142").Value = "EIR"
$ws.Range("E142").Value = "Alt2"
$ws.Range("F142").Value = "`"EIR runs\Alt2 (s28) runs\Alt2_v1`""
$ws.Range("G142").Value = "run374"
$ws.Range("H142").Value = "current"

# Row 143: 2050 EIR Alt2_02 (new anticipated current run)
$ws.Range("A143").Value = "RTP2021"
$ws.Range("B143").Value = 2050
$ws.Range("C143").Value = "2050_TM152_EIR_Alt2_02"
$ws.Range("D143").Value = "EIR"
$ws.Range("E143").Value = "Alt2"
$ws.Range("F143").Value = "`"EIR runs\Alt2 (s28) runs\Alt2_v1`""
$ws.Range("G143").Value = "run374"
$ws.Range("H143").Value = "current"

# --- Move the selection to where the next entry would be typed ---
$ws.Range("C144").Select()
